# Update "flussi del 29 Gennaio 2020" for AttiviPartnerPSP sheet.
# Overwrites Denominazione / CodiceFiscale / NumeroEntiAttivi for every data
# row (2-82), including three brand-new rows appended at the bottom
# (80-82). CodiceFiscale / NumeroEntiAttivi are written with a leading
# quote-prefix so numeric-looking text (e.g. leading zeros) is preserved
# as text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Cells.Item(2, 1).Value2 = 'Ministero dell''Istruzione, Ministero dell''Universita'' e della Ricerca'
$ws.Cells.Item(2, 2).Value2 = '''80185250588'
$ws.Cells.Item(2, 3).Value2 = '''8603'

$ws.Cells.Item(3, 1).Value2 = 'Unione Italiana delle Camere di Commercio Industria, Artigianato e Agricoltura'
$ws.Cells.Item(3, 2).Value2 = '''01484460587'
$ws.Cells.Item(3, 3).Value2 = '''1084'

$ws.Cells.Item(4, 1).Value2 = 'Poste Italiane S.P.A.'
$ws.Cells.Item(4, 2).Value2 = '''97103880585'
$ws.Cells.Item(4, 3).Value2 = '''812'

$ws.Cells.Item(5, 1).Value2 = 'E-Fil S.r.l.'
$ws.Cells.Item(5, 2).Value2 = '''03789750100'
$ws.Cells.Item(5, 3).Value2 = '''633'

$ws.Cells.Item(6, 1).Value2 = 'Regione Lombardia'
$ws.Cells.Item(6, 2).Value2 = '''80050050154'
$ws.Cells.Item(6, 3).Value2 = '''513'

$ws.Cells.Item(7, 1).Value2 = 'Maggioli SPA'
$ws.Cells.Item(7, 2).Value2 = '''06188330150'
$ws.Cells.Item(7, 3).Value2 = '''452'

$ws.Cells.Item(8, 1).Value2 = 'Regione del Veneto'
$ws.Cells.Item(8, 2).Value2 = '''80007580279'
$ws.Cells.Item(8, 3).Value2 = '''430'

$ws.Cells.Item(9, 1).Value2 = 'CREDEMTEL SpA'
$ws.Cells.Item(9, 2).Value2 = '''01378570350'
$ws.Cells.Item(9, 3).Value2 = '''421'

$ws.Cells.Item(10, 1).Value2 = 'Lepida Spa'
$ws.Cells.Item(10, 2).Value2 = '''02770891204'
$ws.Cells.Item(10, 3).Value2 = '''380'

$ws.Cells.Item(11, 1).Value2 = 'Banca Popolare di Sondrio, Società Cooperativa per Azioni'
$ws.Cells.Item(11, 2).Value2 = '''00053810149'
$ws.Cells.Item(11, 3).Value2 = '''369'

$ws.Cells.Item(12, 1).Value2 = 'Provincia Autonoma di Trento'
$ws.Cells.Item(12, 2).Value2 = '''00337460224'
$ws.Cells.Item(12, 3).Value2 = '''336'

$ws.Cells.Item(13, 1).Value2 = 'Dedagroup Public Services S.R.L.'
$ws.Cells.Item(13, 2).Value2 = '''01727860221'
$ws.Cells.Item(13, 3).Value2 = '''255'

$ws.Cells.Item(14, 1).Value2 = 'APKAPPA S.R.L.'
$ws.Cells.Item(14, 2).Value2 = '''08543640158'
$ws.Cells.Item(14, 3).Value2 = '''254'

$ws.Cells.Item(15, 1).Value2 = 'Regione Piemonte'
$ws.Cells.Item(15, 2).Value2 = '''80087670016'
$ws.Cells.Item(15, 3).Value2 = '''238'

$ws.Cells.Item(16, 1).Value2 = 'Progetti e Soluzioni SPA'
$ws.Cells.Item(16, 2).Value2 = '''06423240727'
$ws.Cells.Item(16, 3).Value2 = '''205'

$ws.Cells.Item(17, 1).Value2 = 'Regione Marche'
$ws.Cells.Item(17, 2).Value2 = '''80008630420'
$ws.Cells.Item(17, 3).Value2 = '''184'

$ws.Cells.Item(18, 1).Value2 = 'Alto Adige Riscossioni Spa'
$ws.Cells.Item(18, 2).Value2 = '''02805390214'
$ws.Cells.Item(18, 3).Value2 = '''183'

$ws.Cells.Item(19, 1).Value2 = 'P.A. Digitale spa'
$ws.Cells.Item(19, 2).Value2 = '''06628860964'
$ws.Cells.Item(19, 3).Value2 = '''174'

$ws.Cells.Item(20, 1).Value2 = 'Advanced Systems S.p.A.'
$ws.Cells.Item(20, 2).Value2 = '''03383350638'
$ws.Cells.Item(20, 3).Value2 = '''173'

$ws.Cells.Item(21, 1).Value2 = 'Regione Autonoma della Sardegna'
$ws.Cells.Item(21, 2).Value2 = '''80002870923'
$ws.Cells.Item(21, 3).Value2 = '''155'

$ws.Cells.Item(22, 1).Value2 = 'Siscom SPA'
$ws.Cells.Item(22, 2).Value2 = '''01778000040'
$ws.Cells.Item(22, 3).Value2 = '''128'

$ws.Cells.Item(23, 1).Value2 = 'Regione Puglia'
$ws.Cells.Item(23, 2).Value2 = '''80017210727'
$ws.Cells.Item(23, 3).Value2 = '''127'

$ws.Cells.Item(24, 1).Value2 = 'Next Step Solution'
$ws.Cells.Item(24, 2).Value2 = '''02554480349'
$ws.Cells.Item(24, 3).Value2 = '''118'

$ws.Cells.Item(25, 1).Value2 = 'Regione Toscana'
$ws.Cells.Item(25, 2).Value2 = '''01386030488'
$ws.Cells.Item(25, 3).Value2 = '''115'

$ws.Cells.Item(26, 1).Value2 = 'Regione Basilicata'
$ws.Cells.Item(26, 2).Value2 = '''80002950766'
$ws.Cells.Item(26, 3).Value2 = '''106'

$ws.Cells.Item(27, 1).Value2 = 'Regione Autonoma Friuli-Venezia Giulia'
$ws.Cells.Item(27, 2).Value2 = '''80014930327'
$ws.Cells.Item(27, 3).Value2 = '''98'

$ws.Cells.Item(28, 1).Value2 = 'ANCITEL'
$ws.Cells.Item(28, 2).Value2 = '''07196850585'
$ws.Cells.Item(28, 3).Value2 = '''87'

$ws.Cells.Item(29, 1).Value2 = 'PMPay s.r.l.'
$ws.Cells.Item(29, 2).Value2 = '''08747230962'
$ws.Cells.Item(29, 3).Value2 = '''84'

$ws.Cells.Item(30, 1).Value2 = 'Regione Autonoma Valle D''Aosta'
$ws.Cells.Item(30, 2).Value2 = '''80002270074'
$ws.Cells.Item(30, 3).Value2 = '''81'

$ws.Cells.Item(31, 1).Value2 = 'Intesa Sanpaolo SPA'
$ws.Cells.Item(31, 2).Value2 = '''00799960158'
$ws.Cells.Item(31, 3).Value2 = '''75'

$ws.Cells.Item(32, 1).Value2 = 'ROMA CAPITALE'
$ws.Cells.Item(32, 2).Value2 = '''02438750586'
$ws.Cells.Item(32, 3).Value2 = '''63'

$ws.Cells.Item(33, 1).Value2 = 'Regione Umbria'
$ws.Cells.Item(33, 2).Value2 = '''80000130544'
$ws.Cells.Item(33, 3).Value2 = '''57'

$ws.Cells.Item(34, 1).Value2 = 'NORDCOM'
$ws.Cells.Item(34, 2).Value2 = '''13384100155'
$ws.Cells.Item(34, 3).Value2 = '''56'

$ws.Cells.Item(35, 1).Value2 = 'DCS SOFTWARE E SERVIZI S.R.L.'
$ws.Cells.Item(35, 2).Value2 = '''08063140019'
$ws.Cells.Item(35, 3).Value2 = '''53'

$ws.Cells.Item(36, 1).Value2 = 'Italriscossioni Società Italiana di Fiscalità Locale S.r.l.'
$ws.Cells.Item(36, 2).Value2 = '''06092371001'
$ws.Cells.Item(36, 3).Value2 = '''51'

$ws.Cells.Item(37, 1).Value2 = 'Bluenext S.r.l.'
$ws.Cells.Item(37, 2).Value2 = '''04228480408'
$ws.Cells.Item(37, 3).Value2 = '''46'

$ws.Cells.Item(38, 1).Value2 = 'CINECA consorzio universitario'
$ws.Cells.Item(38, 2).Value2 = '''00317740371'
$ws.Cells.Item(38, 3).Value2 = '''44'

$ws.Cells.Item(39, 1).Value2 = 'UNIMATICA S.P.A'
$ws.Cells.Item(39, 2).Value2 = '''02098391200'
$ws.Cells.Item(39, 3).Value2 = '''43'

$ws.Cells.Item(40, 1).Value2 = 'Consorzio I.T. Srl'
$ws.Cells.Item(40, 2).Value2 = '''01321400192'
$ws.Cells.Item(40, 3).Value2 = '''40'

$ws.Cells.Item(41, 1).Value2 = 'SI.net Servizi Informatici S.r.L.'
$ws.Cells.Item(41, 2).Value2 = '''02743730125'
$ws.Cells.Item(41, 3).Value2 = '''39'

$ws.Cells.Item(42, 1).Value2 = 'Regione Liguria'
$ws.Cells.Item(42, 2).Value2 = '''00849050109'
$ws.Cells.Item(42, 3).Value2 = '''35'

$ws.Cells.Item(43, 1).Value2 = 'Unicredit, Societa'' per Azioni'
$ws.Cells.Item(43, 2).Value2 = '''00348170101'
$ws.Cells.Item(43, 3).Value2 = '''35'

$ws.Cells.Item(44, 1).Value2 = 'Numera Sistemi e Informatica SpA'
$ws.Cells.Item(44, 2).Value2 = '''01265230902'
$ws.Cells.Item(44, 3).Value2 = '''34'

$ws.Cells.Item(45, 1).Value2 = 'Novares Spa'
$ws.Cells.Item(45, 2).Value2 = '''12105121003'
$ws.Cells.Item(45, 3).Value2 = '''31'

$ws.Cells.Item(46, 1).Value2 = 'Comune di Palermo'
$ws.Cells.Item(46, 2).Value2 = '''80016350821'
$ws.Cells.Item(46, 3).Value2 = '''24'

$ws.Cells.Item(47, 1).Value2 = 'Nexi SpA'
$ws.Cells.Item(47, 2).Value2 = '''13212880150'
$ws.Cells.Item(47, 3).Value2 = '''19'

$ws.Cells.Item(48, 1).Value2 = 'Citta'' Metropolitana di Roma Capitale'
$ws.Cells.Item(48, 2).Value2 = '''80034390585'
$ws.Cells.Item(48, 3).Value2 = '''19'

$ws.Cells.Item(49, 1).Value2 = 'ANDREANI TRIBUTI srl'
$ws.Cells.Item(49, 2).Value2 = '''01412920439'
$ws.Cells.Item(49, 3).Value2 = '''18'

$ws.Cells.Item(50, 1).Value2 = 'Regione Lazio'
$ws.Cells.Item(50, 2).Value2 = '''80143490581'
$ws.Cells.Item(50, 3).Value2 = '''18'

$ws.Cells.Item(51, 1).Value2 = 'Si.Form Consulting srl'
$ws.Cells.Item(51, 2).Value2 = '''03943960827'
$ws.Cells.Item(51, 3).Value2 = '''15'

$ws.Cells.Item(52, 1).Value2 = 'Servizi Locali SpA'
$ws.Cells.Item(52, 2).Value2 = '''03170580751'
$ws.Cells.Item(52, 3).Value2 = '''15'

$ws.Cells.Item(53, 1).Value2 = 'Aric Agenzia Regionale di Informatica e Committenza'
$ws.Cells.Item(53, 2).Value2 = '''91022630676'
$ws.Cells.Item(53, 3).Value2 = '''13'

$ws.Cells.Item(54, 1).Value2 = 'Crédit Agricole Group Solutions Società Consortile per azioni'
$ws.Cells.Item(54, 2).Value2 = '''02771790348'
$ws.Cells.Item(54, 3).Value2 = '''13'

$ws.Cells.Item(55, 1).Value2 = 'UBI Banca'
$ws.Cells.Item(55, 2).Value2 = '''03053920165'
$ws.Cells.Item(55, 3).Value2 = '''10'

$ws.Cells.Item(56, 1).Value2 = 'Comune di Catania'
$ws.Cells.Item(56, 2).Value2 = '''00137020871'
$ws.Cells.Item(56, 3).Value2 = '''9'

$ws.Cells.Item(57, 1).Value2 = 'ARCA Servizi s.r.l'
$ws.Cells.Item(57, 2).Value2 = '''09106071005'
$ws.Cells.Item(57, 3).Value2 = '''8'

$ws.Cells.Item(58, 1).Value2 = 'Argentea S.r.l.'
$ws.Cells.Item(58, 2).Value2 = '''02260390220'
$ws.Cells.Item(58, 3).Value2 = '''7'

$ws.Cells.Item(59, 1).Value2 = 'Noviservice srl'
$ws.Cells.Item(59, 2).Value2 = '''02789990922'
$ws.Cells.Item(59, 3).Value2 = '''7'

$ws.Cells.Item(60, 1).Value2 = 'Be Smart s.r.l.'
$ws.Cells.Item(60, 2).Value2 = '''05817461006'
$ws.Cells.Item(60, 3).Value2 = '''7'

$ws.Cells.Item(61, 1).Value2 = 'ARGO SOFTWARE SRL'
$ws.Cells.Item(61, 2).Value2 = '''00838520880'
$ws.Cells.Item(61, 3).Value2 = '''5'

$ws.Cells.Item(62, 1).Value2 = 'CityPoste Payment Digital S.r.l.'
$ws.Cells.Item(62, 2).Value2 = '''02003750672'
$ws.Cells.Item(62, 3).Value2 = '''4'

$ws.Cells.Item(63, 1).Value2 = 'Phoenix IT Solutions S.r.L'
$ws.Cells.Item(63, 2).Value2 = '''07623321218'
$ws.Cells.Item(63, 3).Value2 = '''4'

$ws.Cells.Item(64, 1).Value2 = 'e-SED Società Cooperativa'
$ws.Cells.Item(64, 2).Value2 = '''02695640421'
$ws.Cells.Item(64, 3).Value2 = '''3'

$ws.Cells.Item(65, 1).Value2 = 'ISWEB S.p.A.'
$ws.Cells.Item(65, 2).Value2 = '''01722270665'
$ws.Cells.Item(65, 3).Value2 = '''3'

$ws.Cells.Item(66, 1).Value2 = 'Linea Comune Spa'
$ws.Cells.Item(66, 2).Value2 = '''05591950489'
$ws.Cells.Item(66, 3).Value2 = '''3'

$ws.Cells.Item(67, 1).Value2 = 'KOINE'' SRL'
$ws.Cells.Item(67, 2).Value2 = '''01934790971'
$ws.Cells.Item(67, 3).Value2 = '''2'

$ws.Cells.Item(68, 1).Value2 = 'ICCREA Banca SpA'
$ws.Cells.Item(68, 2).Value2 = '''04774801007'
$ws.Cells.Item(68, 3).Value2 = '''2'

$ws.Cells.Item(69, 1).Value2 = 'Regione Calabria'
$ws.Cells.Item(69, 2).Value2 = '''02205340793'
$ws.Cells.Item(69, 3).Value2 = '''2'

$ws.Cells.Item(70, 1).Value2 = 'Softline srl'
$ws.Cells.Item(70, 2).Value2 = '''12299030150'
$ws.Cells.Item(70, 3).Value2 = '''2'

$ws.Cells.Item(71, 1).Value2 = 'San Marco SPA'
$ws.Cells.Item(71, 2).Value2 = '''04142440728'
$ws.Cells.Item(71, 3).Value2 = '''1'

$ws.Cells.Item(72, 1).Value2 = 'Engineering Ingegneria Informatica SpA'
$ws.Cells.Item(72, 2).Value2 = '''00967720285'
$ws.Cells.Item(72, 3).Value2 = '''1'

$ws.Cells.Item(73, 1).Value2 = 'Società Almaviva S.p.A.'
$ws.Cells.Item(73, 2).Value2 = '''08450891000'
$ws.Cells.Item(73, 3).Value2 = '''1'

$ws.Cells.Item(74, 1).Value2 = 'Agenzia Italiana del Farmaco - AIFA'
$ws.Cells.Item(74, 2).Value2 = '''97345810580'
$ws.Cells.Item(74, 3).Value2 = '''1'

$ws.Cells.Item(75, 1).Value2 = 'I.C.A. - Imposte Comunali Affini – s.r.l.'
$ws.Cells.Item(75, 2).Value2 = '''02478610583'
$ws.Cells.Item(75, 3).Value2 = '''1'

$ws.Cells.Item(76, 1).Value2 = 'BANCA MONTE DEI PASCHI DI SIENA'
$ws.Cells.Item(76, 2).Value2 = '''00884060526'
$ws.Cells.Item(76, 3).Value2 = '''1'

$ws.Cells.Item(77, 1).Value2 = 'Open Software S.r.l.'
$ws.Cells.Item(77, 2).Value2 = '''02810000279'
$ws.Cells.Item(77, 3).Value2 = '''1'

$ws.Cells.Item(78, 1).Value2 = 'Banco BPM Società per Azioni'
$ws.Cells.Item(78, 2).Value2 = '''09722490969'
$ws.Cells.Item(78, 3).Value2 = '''1'

$ws.Cells.Item(79, 1).Value2 = 'WAN S.r.l.'
$ws.Cells.Item(79, 2).Value2 = '''03805290040'
$ws.Cells.Item(79, 3).Value2 = '''1'

$ws.Cells.Item(80, 1).Value2 = 'Banca Nazionale del Lavoro S.p.A.'
$ws.Cells.Item(80, 2).Value2 = '''09339391006'
$ws.Cells.Item(80, 3).Value2 = '''1'

$ws.Cells.Item(81, 1).Value2 = 'Ministero dello Sviluppo Economico'
$ws.Cells.Item(81, 2).Value2 = '''80230390587'
$ws.Cells.Item(81, 3).Value2 = '''1'

$ws.Cells.Item(82, 1).Value2 = 'MegASP S.r.l.'
$ws.Cells.Item(82, 2).Value2 = '''09898030151'
$ws.Cells.Item(82, 3).Value2 = '''1'
